# Fill in the previously-empty "class" column (G) cells with "lipid/free"
# for the glycan rows that represent free/lipid-linked glycans (rows whose
# G cell was an empty inline string in the original workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4, 5, 7, 8, 10, 11, 14, 16, 17, 21, 26, 27, 31)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "lipid/free"
}
